# Updates market-derived Leve profit columns (H-N) across all 8 Job sheets.
# Values correspond to a scheduled data refresh (Universalis price snapshot).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1104.6316   # H15
$ws.Cells.Item(15, 9).Value = 1104.6316   # I15
$ws.Cells.Item(15, 11).Value = 3313.8948   # K15
$ws.Cells.Item(15, 13).Value = -3144.8948   # M15
$ws.Cells.Item(43, 8).Value = 4000.25   # H43
$ws.Cells.Item(43, 9).Value = 4000.25   # I43
$ws.Cells.Item(43, 10).Value = 0   # J43
$ws.Cells.Item(43, 11).Value = 4000.25   # K43
$ws.Cells.Item(43, 12).Value = 0   # L43
$ws.Cells.Item(43, 13).Value = -3931.25   # M43
$ws.Cells.Item(43, 14).Value = $null   # N43 (cleared)
$ws.Cells.Item(92, 8).Value = 967.53125   # H92
$ws.Cells.Item(92, 9).Value = 989.86957   # I92
$ws.Cells.Item(92, 10).Value = 910.44446   # J92
$ws.Cells.Item(92, 11).Value = 989.86957   # K92
$ws.Cells.Item(92, 12).Value = 910.44446   # L92
$ws.Cells.Item(92, 13).Value = 258.13043   # M92
$ws.Cells.Item(92, 14).Value = -3406.44446   # N92
$ws.Cells.Item(113, 8).Value = 126770.375   # H113
$ws.Cells.Item(113, 9).Value = 334727.66   # I113
$ws.Cells.Item(113, 10).Value = 1996   # J113
$ws.Cells.Item(113, 11).Value = 334727.66   # K113
$ws.Cells.Item(113, 12).Value = 1996   # L113
$ws.Cells.Item(113, 13).Value = -331473.66   # M113
$ws.Cells.Item(113, 14).Value = -8504   # N113
$ws.Cells.Item(137, 8).Value = 954.6405999999999   # H137
$ws.Cells.Item(137, 9).Value = 952.1177   # I137
$ws.Cells.Item(137, 11).Value = 2856.3531   # K137
$ws.Cells.Item(137, 13).Value = -306.3531000000003   # M137
$ws.Cells.Item(138, 8).Value = 1556.8833   # H138
$ws.Cells.Item(138, 9).Value = 987.4186   # I138
$ws.Cells.Item(138, 10).Value = 2997.2942   # J138
$ws.Cells.Item(138, 11).Value = 2962.2558   # K138
$ws.Cells.Item(138, 12).Value = 8991.882599999999   # L138
$ws.Cells.Item(138, 13).Value = 2177.7442   # M138
$ws.Cells.Item(138, 14).Value = -19271.8826   # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 513.83   # H32
$ws.Cells.Item(32, 9).Value = 514.2258   # I32
$ws.Cells.Item(32, 10).Value = 508.57144   # J32
$ws.Cells.Item(32, 11).Value = 514.2258   # K32
$ws.Cells.Item(32, 12).Value = 508.57144   # L32
$ws.Cells.Item(32, 13).Value = -227.2258   # M32
$ws.Cells.Item(32, 14).Value = -1082.57144   # N32
$ws.Cells.Item(64, 8).Value = 40830   # H64
$ws.Cells.Item(64, 10).Value = 40830   # J64
$ws.Cells.Item(64, 12).Value = 40830   # L64
$ws.Cells.Item(64, 14).Value = -41326   # N64
$ws.Cells.Item(67, 8).Value = 40830   # H67
$ws.Cells.Item(67, 10).Value = 40830   # J67
$ws.Cells.Item(67, 12).Value = 40830   # L67
$ws.Cells.Item(67, 14).Value = -42546   # N67
$ws.Cells.Item(74, 8).Value = 563.6   # H74
$ws.Cells.Item(74, 9).Value = 432.6842   # I74
$ws.Cells.Item(74, 10).Value = 789.7273   # J74
$ws.Cells.Item(74, 11).Value = 432.6842   # K74
$ws.Cells.Item(74, 12).Value = 789.7273   # L74
$ws.Cells.Item(74, 13).Value = 441.3158   # M74
$ws.Cells.Item(74, 14).Value = -2537.7273   # N74
$ws.Cells.Item(77, 8).Value = 563.6   # H77
$ws.Cells.Item(77, 9).Value = 432.6842   # I77
$ws.Cells.Item(77, 10).Value = 789.7273   # J77
$ws.Cells.Item(77, 11).Value = 2163.421   # K77
$ws.Cells.Item(77, 12).Value = 3948.6365   # L77
$ws.Cells.Item(77, 13).Value = 2204.579   # M77
$ws.Cells.Item(77, 14).Value = -12684.6365   # N77

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(54, 8).Value = 3784.7   # H54
$ws.Cells.Item(54, 9).Value = 2427.4443   # I54
$ws.Cells.Item(54, 11).Value = 2427.4443   # K54
$ws.Cells.Item(54, 13).Value = -1943.4443   # M54
$ws.Cells.Item(62, 8).Value = 48000   # H62
$ws.Cells.Item(62, 10).Value = 48000   # J62
$ws.Cells.Item(62, 12).Value = 48000   # L62
$ws.Cells.Item(62, 14).Value = -49372   # N62
$ws.Cells.Item(65, 8).Value = 48000   # H65
$ws.Cells.Item(65, 10).Value = 48000   # J65
$ws.Cells.Item(65, 12).Value = 144000   # L65
$ws.Cells.Item(65, 14).Value = -150864   # N65

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 874.50726   # H58
$ws.Cells.Item(58, 9).Value = 770.36365   # I58
$ws.Cells.Item(58, 10).Value = 1283.6428   # J58
$ws.Cells.Item(58, 11).Value = 770.36365   # K58
$ws.Cells.Item(58, 12).Value = 1283.6428   # L58
$ws.Cells.Item(58, 13).Value = -567.36365   # M58
$ws.Cells.Item(58, 14).Value = -1689.6428   # N58
$ws.Cells.Item(134, 8).Value = 954.3542   # H134
$ws.Cells.Item(134, 9).Value = 926.0714   # I134
$ws.Cells.Item(134, 10).Value = 1152.3334   # J134
$ws.Cells.Item(134, 11).Value = 2778.2142   # K134
$ws.Cells.Item(134, 12).Value = 3457.0002   # L134
$ws.Cells.Item(134, 13).Value = -243.2142000000003   # M134
$ws.Cells.Item(134, 14).Value = -8527.0002   # N134
$ws.Cells.Item(136, 8).Value = 874.50726   # H136
$ws.Cells.Item(136, 9).Value = 770.36365   # I136
$ws.Cells.Item(136, 10).Value = 1283.6428   # J136
$ws.Cells.Item(136, 11).Value = 2311.09095   # K136
$ws.Cells.Item(136, 12).Value = 3850.9284   # L136
$ws.Cells.Item(136, 13).Value = 238.9090500000002   # M136
$ws.Cells.Item(136, 14).Value = -8950.928400000001   # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 952.4878   # H5
$ws.Cells.Item(5, 9).Value = 481.58334   # I5
$ws.Cells.Item(5, 10).Value = 1617.2941   # J5
$ws.Cells.Item(5, 11).Value = 1444.75002   # K5
$ws.Cells.Item(5, 12).Value = 4851.8823   # L5
$ws.Cells.Item(5, 13).Value = -1332.75002   # M5
$ws.Cells.Item(5, 14).Value = -5075.8823   # N5
$ws.Cells.Item(122, 8).Value = 493.57144   # H122
$ws.Cells.Item(122, 9).Value = 491   # I122
$ws.Cells.Item(122, 10).Value = 500   # J122
$ws.Cells.Item(122, 11).Value = 4419   # K122
$ws.Cells.Item(122, 12).Value = 4500   # L122
$ws.Cells.Item(122, 13).Value = -1969   # M122
$ws.Cells.Item(122, 14).Value = -9400   # N122
$ws.Cells.Item(129, 8).Value = 1805.8334   # H129
$ws.Cells.Item(129, 9).Value = 652.8570999999999   # I129
$ws.Cells.Item(129, 10).Value = 3420   # J129
$ws.Cells.Item(129, 11).Value = 1958.5713   # K129
$ws.Cells.Item(129, 12).Value = 10260   # L129
$ws.Cells.Item(129, 13).Value = 3041.4287   # M129
$ws.Cells.Item(129, 14).Value = -20260   # N129
$ws.Cells.Item(131, 8).Value = 1318.1086   # H131
$ws.Cells.Item(131, 10).Value = 1364.9642   # J131
$ws.Cells.Item(131, 12).Value = 4094.8926   # L131
$ws.Cells.Item(131, 14).Value = -14174.8926   # N131
$ws.Cells.Item(133, 8).Value = 4778.3   # H133
$ws.Cells.Item(133, 9).Value = 3261   # I133
$ws.Cells.Item(133, 10).Value = 5428.5713   # J133
$ws.Cells.Item(133, 11).Value = 9783   # K133
$ws.Cells.Item(133, 12).Value = 16285.7139   # L133
$ws.Cells.Item(133, 13).Value = -4723   # M133
$ws.Cells.Item(133, 14).Value = -26405.7139   # N133
$ws.Cells.Item(135, 8).Value = 952.4878   # H135
$ws.Cells.Item(135, 9).Value = 481.58334   # I135
$ws.Cells.Item(135, 10).Value = 1617.2941   # J135
$ws.Cells.Item(135, 11).Value = 4334.25006   # K135
$ws.Cells.Item(135, 12).Value = 14555.6469   # L135
$ws.Cells.Item(135, 13).Value = -1799.25006   # M135
$ws.Cells.Item(135, 14).Value = -19625.6469   # N135
$ws.Cells.Item(140, 8).Value = 4025.6216   # H140
$ws.Cells.Item(140, 9).Value = 4894.76   # I140
$ws.Cells.Item(140, 10).Value = 2214.9167   # J140
$ws.Cells.Item(140, 11).Value = 14684.28   # K140
$ws.Cells.Item(140, 12).Value = 6644.750100000001   # L140
$ws.Cells.Item(140, 13).Value = -9504.280000000001   # M140
$ws.Cells.Item(140, 14).Value = -17004.7501   # N140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2259.7778   # H102
$ws.Cells.Item(102, 9).Value = 1778.2858   # I102
$ws.Cells.Item(102, 10).Value = 2778.3076   # J102
$ws.Cells.Item(102, 11).Value = 1778.2858   # K102
$ws.Cells.Item(102, 12).Value = 2778.3076   # L102
$ws.Cells.Item(102, 13).Value = -156.2858000000001   # M102
$ws.Cells.Item(102, 14).Value = -6022.3076   # N102
$ws.Cells.Item(132, 8).Value = 1705.2898   # H132
$ws.Cells.Item(132, 9).Value = 1646.4395   # I132
$ws.Cells.Item(132, 11).Value = 4939.318499999999   # K132
$ws.Cells.Item(132, 13).Value = -2409.318499999999   # M132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(45, 8).Value = 6233.4   # H45
$ws.Cells.Item(45, 9).Value = 3530.5   # I45
$ws.Cells.Item(45, 10).Value = 6909.125   # J45
$ws.Cells.Item(45, 11).Value = 3530.5   # K45
$ws.Cells.Item(45, 12).Value = 6909.125   # L45
$ws.Cells.Item(45, 13).Value = -3123.5   # M45
$ws.Cells.Item(45, 14).Value = -7723.125   # N45
$ws.Cells.Item(61, 8).Value = 1457.8096   # H61
$ws.Cells.Item(61, 9).Value = 1340.6   # I61
$ws.Cells.Item(61, 10).Value = 1750.8334   # J61
$ws.Cells.Item(61, 11).Value = 1340.6   # K61
$ws.Cells.Item(61, 12).Value = 1750.8334   # L61
$ws.Cells.Item(61, 13).Value = -1138.6   # M61
$ws.Cells.Item(61, 14).Value = -2154.8334   # N61
$ws.Cells.Item(113, 8).Value = 1457.8096   # H113
$ws.Cells.Item(113, 9).Value = 1340.6   # I113
$ws.Cells.Item(113, 10).Value = 1750.8334   # J113
$ws.Cells.Item(113, 11).Value = 1340.6   # K113
$ws.Cells.Item(113, 12).Value = 1750.8334   # L113
$ws.Cells.Item(113, 13).Value = 829.4000000000001   # M113
$ws.Cells.Item(113, 14).Value = -6090.8334   # N113
$ws.Cells.Item(132, 8).Value = 1945.5536   # H132
$ws.Cells.Item(132, 9).Value = 1933.0377   # I132
$ws.Cells.Item(132, 11).Value = 5799.1131   # K132
$ws.Cells.Item(132, 13).Value = -3269.1131   # M132
$ws.Cells.Item(136, 8).Value = 1473.0344   # H136
$ws.Cells.Item(136, 9).Value = 1339.1538   # I136
$ws.Cells.Item(136, 10).Value = 2633.3333   # J136
$ws.Cells.Item(136, 11).Value = 4017.4614   # K136
$ws.Cells.Item(136, 12).Value = 7899.999899999999   # L136
$ws.Cells.Item(136, 13).Value = -1467.4614   # M136
$ws.Cells.Item(136, 14).Value = -12999.9999   # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 400680   # H81
$ws.Cells.Item(81, 9).Value = 250850.25   # I81
$ws.Cells.Item(81, 10).Value = 999999   # J81
$ws.Cells.Item(81, 11).Value = 501700.5   # K81
$ws.Cells.Item(81, 12).Value = 1999998   # L81
$ws.Cells.Item(81, 13).Value = -500639.5   # M81
$ws.Cells.Item(81, 14).Value = -2002120   # N81
$ws.Cells.Item(84, 8).Value = 400680   # H84
$ws.Cells.Item(84, 9).Value = 250850.25   # I84
$ws.Cells.Item(84, 10).Value = 999999   # J84
$ws.Cells.Item(84, 11).Value = 2508502.5   # K84
$ws.Cells.Item(84, 12).Value = 9999990   # L84
$ws.Cells.Item(84, 13).Value = -2503198.5   # M84
$ws.Cells.Item(84, 14).Value = -10010598   # N84

Write-Output "Applied 201 cell updates"
